$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 492.5
$ws.Range("I9").Value = 700
$ws.Range("J9").Value = 285
$ws.Range("K9").Value = 700
$ws.Range("L9").Value = 285
$ws.Range("M9").Value = -531
$ws.Range("N9").Value = -623

$ws.Range("H17").Value = 1788.75
$ws.Range("J17").Value = 1788.75
$ws.Range("L17").Value = 5366.25
$ws.Range("N17").Value = -5702.25

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H33").Value = 262.42856
$ws.Range("I33").Value = 262.42856
$ws.Range("K33").Value = 262.42856
$ws.Range("M33").Value = -33.42856

$ws.Range("H53").Value = 177
$ws.Range("J53").Value = 226.33333
$ws.Range("L53").Value = 226.33333
$ws.Range("N53").Value = -1500.33333

$ws.Range("H70").Value = 11398
$ws.Range("I70").Value = 6748.75
$ws.Range("K70").Value = 20246.25
$ws.Range("M70").Value = -19976.25

$ws.Range("H73").Value = 11398
$ws.Range("I73").Value = 6748.75
$ws.Range("K73").Value = 20246.25
$ws.Range("M73").Value = -19310.25

$ws.Range("H138").Value = 3703.6
$ws.Range("I138").Value = 3155
$ws.Range("K138").Value = 9465
$ws.Range("M138").Value = -4325


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 975
$ws.Range("I12").Value = 450
$ws.Range("J12").Value = 1500
$ws.Range("K12").Value = 450
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = -277
$ws.Range("N12").Value = -1846

$ws.Range("H102").Value = 3410.3333
$ws.Range("I102").Value = 3410.3333
$ws.Range("K102").Value = 3410.3333
$ws.Range("M102").Value = -1788.3333

$ws.Range("H109").Value = 23995
$ws.Range("J109").Value = 23995
$ws.Range("L109").Value = 23995
$ws.Range("N109").Value = -26769

$ws.Range("H122").Value = 9815.091
$ws.Range("J122").Value = 9000
$ws.Range("L122").Value = 27000
$ws.Range("N122").Value = -31900


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 129
$ws.Range("I5").Value = 129
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 129
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -16
$ws.Range("N5").ClearContents()

$ws.Range("H11").Value = 337.66666
$ws.Range("J11").Value = 478
$ws.Range("L11").Value = 478
$ws.Range("N11").Value = -758

$ws.Range("H134").Value = 9310.625
$ws.Range("I134").Value = 2247.5
$ws.Range("K134").Value = 6742.5
$ws.Range("M134").Value = -4207.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1026
$ws.Range("I8").Value = 1026
$ws.Range("K8").Value = 1026
$ws.Range("M8").Value = -886

$ws.Range("H16").Value = 583.3333
$ws.Range("I16").Value = 583.3333
$ws.Range("K16").Value = 583.3333
$ws.Range("M16").Value = -296.3333

$ws.Range("H19").Value = 345
$ws.Range("I19").Value = 451.5
$ws.Range("J19").Value = 25.5
$ws.Range("K19").Value = 451.5
$ws.Range("L19").Value = 25.5
$ws.Range("M19").Value = -281.5
$ws.Range("N19").Value = -365.5

$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -765

$ws.Range("H24").Value = 345
$ws.Range("I24").Value = 451.5
$ws.Range("J24").Value = 25.5
$ws.Range("K24").Value = 451.5
$ws.Range("L24").Value = 25.5
$ws.Range("M24").Value = -281.5
$ws.Range("N24").Value = -365.5

$ws.Range("H50").Value = 23579.5
$ws.Range("J50").Value = 25088.334
$ws.Range("L50").Value = 25088.334
$ws.Range("N50").Value = -26338.334

$ws.Range("H99").Value = 3000
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996

$ws.Range("H105").Value = 1312.8334
$ws.Range("I105").Value = 1469.5
$ws.Range("J105").Value = 999.5
$ws.Range("K105").Value = 1469.5
$ws.Range("L105").Value = 999.5
$ws.Range("M105").Value = 277.5
$ws.Range("N105").Value = -4493.5

$ws.Range("H113").Value = 583.3333
$ws.Range("I113").Value = 583.3333
$ws.Range("K113").Value = 583.3333
$ws.Range("M113").Value = 1586.6667

$ws.Range("H126").Value = 3000
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 55
$ws.Range("K5").Value = 55
$ws.Range("M5").Value = 57

$ws.Range("H9").Value = 2384.6667
$ws.Range("I9").Value = 2384.6667
$ws.Range("K9").Value = 2384.6667
$ws.Range("M9").Value = -2214.6667

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H126").Value = 3245
$ws.Range("I126").Value = 3245
$ws.Range("K126").Value = 9735
$ws.Range("M126").Value = -7265


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 8000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -7888
$ws.Range("N7").ClearContents()

$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550

$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 8000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21530
$ws.Range("N126").ClearContents()

$ws.Range("H138").Value = 99429
$ws.Range("J138").Value = 99429
$ws.Range("L138").Value = 99429
$ws.Range("N138").Value = -109709


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 31588
$ws.Range("I4").Value = 33633.23
$ws.Range("K4").Value = 33633.23
$ws.Range("M4").Value = -33520.23

$ws.Range("H122").Value = 7599.8
$ws.Range("I122").Value = 6333
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 18999
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -16549
$ws.Range("N122").Value = -33400

